$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 10731.5
$ws.Range("I5").Value = 10411
$ws.Range("K5").Value = 10411
$ws.Range("M5").Value = -10296
$ws.Range("H28").Value = 1851.5
$ws.Range("I28").Value = 1449.2307
$ws.Range("J28").Value = 2432.5557
$ws.Range("K28").Value = 1449.2307
$ws.Range("L28").Value = 2432.5557
$ws.Range("M28").Value = -964.2307000000001
$ws.Range("N28").Value = -3402.5557
$ws.Range("H32").Value = 2499
$ws.Range("I32").Value = 1374.5
$ws.Range("K32").Value = 1374.5
$ws.Range("M32").Value = -1048.5
$ws.Range("H116").Value = 1443301.5
$ws.Range("I116").Value = 1672890.2
$ws.Range("K116").Value = 1672890.2
$ws.Range("M116").Value = -1669448.2
$ws.Range("H133").Value = 120780
$ws.Range("J133").Value = 120780
$ws.Range("L133").Value = 120780
$ws.Range("N133").Value = -130900

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 604669.75
$ws.Range("I2").Value = 875876.7
$ws.Range("J2").Value = 1987.6666
$ws.Range("K2").Value = 875876.7
$ws.Range("L2").Value = 1987.6666
$ws.Range("M2").Value = -875763.7
$ws.Range("N2").Value = -2213.6666
$ws.Range("H45").Value = 4332.5
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H63").Value = 3171.2856
$ws.Range("I63").Value = 1841
$ws.Range("K63").Value = 1841
$ws.Range("M63").Value = -1155
$ws.Range("H66").Value = 3171.2856
$ws.Range("I66").Value = 1841
$ws.Range("K66").Value = 9205
$ws.Range("M66").Value = -5773
$ws.Range("H92").Value = 30000000
$ws.Range("J92").Value = 30000000
$ws.Range("L92").Value = 30000000
$ws.Range("N92").Value = -30004992
$ws.Range("H116").Value = 604669.75
$ws.Range("I116").Value = 875876.7
$ws.Range("J116").Value = 1987.6666
$ws.Range("K116").Value = 875876.7
$ws.Range("L116").Value = 1987.6666
$ws.Range("M116").Value = -873582.7
$ws.Range("N116").Value = -6575.6666
$ws.Range("H132").Value = 13847.852
$ws.Range("I132").Value = 18655.152
$ws.Range("J132").Value = 6293.524
$ws.Range("K132").Value = 55965.45599999999
$ws.Range("L132").Value = 18880.572
$ws.Range("M132").Value = -53435.45599999999
$ws.Range("N132").Value = -23940.572

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 604669.75
$ws.Range("I3").Value = 875876.7
$ws.Range("J3").Value = 1987.6666
$ws.Range("K3").Value = 875876.7
$ws.Range("L3").Value = 1987.6666
$ws.Range("M3").Value = -875762.7
$ws.Range("N3").Value = -2215.6666
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2453.6
$ws.Range("I16").Value = 2379.5
$ws.Range("K16").Value = 2379.5
$ws.Range("M16").Value = -2092.5
$ws.Range("H31").Value = 6525.2686
$ws.Range("I31").Value = 2984.125
$ws.Range("K31").Value = 2984.125
$ws.Range("M31").Value = -2689.125
$ws.Range("H34").Value = 6525.2686
$ws.Range("I34").Value = 2984.125
$ws.Range("K34").Value = 2984.125
$ws.Range("M34").Value = -2782.125
$ws.Range("H58").Value = 517169.78
$ws.Range("I58").Value = 528090.6
$ws.Range("J58").Value = 506794.94
$ws.Range("K58").Value = 528090.6
$ws.Range("L58").Value = 506794.94
$ws.Range("M58").Value = -527887.6
$ws.Range("N58").Value = -507200.94
$ws.Range("H60").Value = 75819.75
$ws.Range("J60").Value = 300000
$ws.Range("L60").Value = 300000
$ws.Range("N60").Value = -301022
$ws.Range("H99").Value = 5508.241
$ws.Range("I99").Value = 3640.8667
$ws.Range("J99").Value = 7509
$ws.Range("K99").Value = 3640.8667
$ws.Range("L99").Value = 7509
$ws.Range("M99").Value = -2142.8667
$ws.Range("N99").Value = -10505
$ws.Range("H113").Value = 2453.6
$ws.Range("I113").Value = 2379.5
$ws.Range("K113").Value = 2379.5
$ws.Range("M113").Value = -209.5
$ws.Range("H122").Value = 3713.8276
$ws.Range("I122").Value = 2259.7778
$ws.Range("J122").Value = 6093.1816
$ws.Range("K122").Value = 6779.3334
$ws.Range("L122").Value = 18279.5448
$ws.Range("M122").Value = -4329.3334
$ws.Range("N122").Value = -23179.5448
$ws.Range("H126").Value = 5508.241
$ws.Range("I126").Value = 3640.8667
$ws.Range("J126").Value = 7509
$ws.Range("K126").Value = 10922.6001
$ws.Range("L126").Value = 22527
$ws.Range("M126").Value = -8452.6001
$ws.Range("N126").Value = -27467
$ws.Range("H134").Value = 1801.5
$ws.Range("I134").Value = 1801.5
$ws.Range("K134").Value = 5404.5
$ws.Range("M134").Value = -2869.5
$ws.Range("H136").Value = 517169.78
$ws.Range("I136").Value = 528090.6
$ws.Range("J136").Value = 506794.94
$ws.Range("K136").Value = 1584271.8
$ws.Range("L136").Value = 1520384.82
$ws.Range("M136").Value = -1581721.8
$ws.Range("N136").Value = -1525484.82

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 24498
$ws.Range("J112").Value = 24498
$ws.Range("L112").Value = 73494
$ws.Range("N112").Value = -75710

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 26325994
$ws.Range("I102").Value = 35724636
$ws.Range("K102").Value = 35724636
$ws.Range("M102").Value = -35723014
$ws.Range("H123").Value = 52372.168
$ws.Range("J123").Value = 52372.168
$ws.Range("L123").Value = 52372.168
$ws.Range("N123").Value = -57272.168

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4104.7026
$ws.Range("I132").Value = 3565.9167
$ws.Range("K132").Value = 10697.7501
$ws.Range("M132").Value = -8167.750100000001
$ws.Range("H136").Value = 3305.3845
$ws.Range("I136").Value = 1944.1052
$ws.Range("J136").Value = 7000.2856
$ws.Range("K136").Value = 5832.3156
$ws.Range("L136").Value = 21000.8568
$ws.Range("M136").Value = -3282.3156
$ws.Range("N136").Value = -26100.8568

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 99999
$ws.Range("J64").Value = 99999
$ws.Range("L64").Value = 99999
$ws.Range("N64").Value = -100495
$ws.Range("H67").Value = 99999
$ws.Range("J67").Value = 99999
$ws.Range("L67").Value = 99999
$ws.Range("N67").Value = -101715
$ws.Range("H107").Value = 3072.7144
$ws.Range("I107").Value = 3253.4211
$ws.Range("K107").Value = 9760.263300000001
$ws.Range("M107").Value = -7840.263300000001
$ws.Range("H113").Value = 1098.9333
$ws.Range("I113").Value = 1097.7778
$ws.Range("J113").Value = 1100.6666
$ws.Range("K113").Value = 3293.3334
$ws.Range("L113").Value = 3301.9998
$ws.Range("M113").Value = -1123.3334
$ws.Range("N113").Value = -7641.9998

